# Applies the cryptos list update described in the commit
# "Updated cryptos list on Mon Feb 19 09:58:20 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'52.340.48"
$ws.Range("E2").Value = "'  +0.81%  "
# Row 3
$ws.Range("D3").Value = "'2.913.18"
$ws.Range("E3").Value = "'  +3.42%  "
# Row 4
$ws.Range("E4").Value = "'  +0.11%  "
# Row 5
$ws.Range("D5").Value = "'351.73"
$ws.Range("E5").Value = "'  -1.46%  "
# Row 6
$ws.Range("D6").Value = "'112.27"
$ws.Range("E6").Value = "'  +1.47%  "
# Row 7
$ws.Range("D7").Value = "'0.562"
$ws.Range("E7").Value = "'  +0.70%  "
# Row 8
$ws.Range("E8").Value = "'  +0.03%  "
# Row 9
$ws.Range("D9").Value = "'0.632"
$ws.Range("E9").Value = "'  -0.52%  "
# Row 10
$ws.Range("D10").Value = "'40.06"
$ws.Range("E10").Value = "'  -0.76%  "
# Row 11
$ws.Range("D11").Value = "'0.0865"
$ws.Range("E11").Value = "'  +2.86%  "
# Row 12
$ws.Range("D12").Value = "'0.135"
$ws.Range("E12").Value = "'  +0.12%  "
# Row 13
$ws.Range("D13").Value = "'19.95"
$ws.Range("E13").Value = "'  -0.44%  "
# Row 14
$ws.Range("D14").Value = "'7.80"
$ws.Range("E14").Value = "'  -0.23%  "
# Row 15
$ws.Range("D15").Value = "'3.374.90"
$ws.Range("E15").Value = "'  +3.62%  "
# Row 16
$ws.Range("E16").Value = "'  +6.79%  "
# Row 17
$ws.Range("D17").Value = "'2.936.65"
$ws.Range("E17").Value = "'  +4.22%  "
# Row 18
$ws.Range("D18").Value = "'52.378.90"
$ws.Range("E18").Value = "'  +0.95%  "
# Row 19
$ws.Range("D19").Value = "'3.35"
$ws.Range("E19").Value = "'  +4.64%  "
# Row 20
$ws.Range("D20").Value = "'7.64"
$ws.Range("E20").Value = "'  -0.71%  "
# Row 21
$ws.Range("D21").Value = "'14.18"
$ws.Range("E21").Value = "'  +3.47%  "
# Row 22
$ws.Range("D22").Value = "'0.0₃0980"
$ws.Range("E22").Value = "'  -0.04%  "
# Row 23
$ws.Range("D23").Value = "'71.12"
# Row 24
$ws.Range("D24").Value = "'270.52"
$ws.Range("E24").Value = "'  +0.67%  "
# Row 25
$ws.Range("D25").Value = "'2.79"
$ws.Range("E25").Value = "'  +0.65%  "
# Row 26
$ws.Range("D26").Value = "'26.76"
$ws.Range("E26").Value = "'  +1.91%  "
# Row 27
$ws.Range("E27").Value = "'  +2.68%  "
# Row 28
$ws.Range("E28").Value = "'  -0.10%  "
# Row 29
$ws.Range("D29").Value = "'10.64"
$ws.Range("E29").Value = "'  +1.85%  "
# Row 30
$ws.Range("B30").Value = "RenderToken"
$ws.Range("C30").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D30").Value = "'6.37"
$ws.Range("E30").Value = "'  +12.45%  "
# Row 31
$ws.Range("B31").Value = "InjectiveProtocol"
$ws.Range("C31").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D31").Value = "'37.71"
$ws.Range("E31").Value = "'  -1.64%  "
# Row 32
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").Value = "'6.59"
$ws.Range("E32").Value = "'  +5.91%  "
# Row 33
$ws.Range("E33").Value = "'  +0.61%  "
# Row 34
$ws.Range("E34").Value = "'  +10.58%  "
# Row 35
$ws.Range("D35").Value = "'53.18"
$ws.Range("E35").Value = "'  +0.99%  "
# Row 36
$ws.Range("E36").Value = "'  +1.28%  "
# Row 37
$ws.Range("E37").Value = "'  -0.04%  "
# Row 38
$ws.Range("D38").Value = "'3.32"
$ws.Range("E38").Value = "'  +5.50%  "
# Row 39
$ws.Range("D39").Value = "'18.73"
$ws.Range("E39").Value = "'  -0.46%  "
# Row 40
$ws.Range("E40").Value = "'  +2.64%  "
# Row 41
$ws.Range("D41").Value = "'2.85"
$ws.Range("E41").Value = "'  +13.25%  "
# Row 42
$ws.Range("D42").Value = "'23.78"
$ws.Range("E42").Value = "'  +8.05%  "
# Row 43
$ws.Range("E43").Value = "'  +1.50%  "
# Row 44
$ws.Range("B44").Value = "ApeXProtocol"
$ws.Range("C44").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D44").Value = "'2.62"
$ws.Range("E44").Value = "'  +7.49%  "
# Row 45
$ws.Range("B45").Value = "Monero"
$ws.Range("C45").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D45").Value = "'121.21"
$ws.Range("E45").Value = "'  +0.60%  "
# Row 46
$ws.Range("E46").Value = "'  -0.25%  "
# Row 47
$ws.Range("D47").Value = "'3.55"
$ws.Range("E47").Value = "'  +4.38%  "
# Row 48
$ws.Range("D48").Value = "'2.198.61"
$ws.Range("E48").Value = "'  +4.08%  "
# Row 49
$ws.Range("D49").Value = "'0.266"
$ws.Range("E49").Value = "'  +23.95%  "
# Row 50
$ws.Range("D50").Value = "'0.0339"
$ws.Range("E50").Value = "'  +12.27%  "
# Row 51
$ws.Range("D51").Value = "'0.964"
$ws.Range("E51").Value = "'  +2.69%  "
